# Generate Report for handback
# Updates the "zh-cn" and "de-de" localization status sheets to reflect that
# a.md.md / b.md.md have now been handed back (in sync with en-US), and
# records the new handback file/datetime information in the previously
# empty "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet ----------------------------------------------------
# The Overview sheet mirrors the same "Ready for handoff" status text (it
# shares the underlying string with the per-language sheets), so it also
# flips to the handed-back status for both a.md.md and b.md.md rows.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# ---- zh-cn sheet -----------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Status column (B) for both tracked files moves from "Ready for handoff"
# to "Handed back: in sync with en-US".
$zh.Range("B2").Value = $statusText
$zh.Range("B3").Value = $statusText

# Newly populated "Latest Target File" (E) / "Latest Handback File" (F)
# columns, mirroring the existing handoff file hyperlinks (A/C) since the
# handback target is the same source/handoff artifact.
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ff8a855593417b12f2c3f4504fdc7a9c0dd4205/e2e/a.md.md", "", "", "a.md.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/19d5eb49cba9812a89c02881854f6b7009bfea6b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/6ff8a855593417b12f2c3f4504fdc7a9c0dd4205/e2e/a.md.md", "", "", "a.md.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/19d5eb49cba9812a89c02881854f6b7009bfea6b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf") | Out-Null

# Latest Handback DateTime (G) now has a real timestamp instead of the
# zero-date placeholder.
$zh.Range("G2").Value = "2016-01-19 04:19:37"
$zh.Range("G3").Value = "2016-01-19 04:19:37"

# ---- de-de sheet -------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $statusText
$de.Range("B3").Value = $statusText

$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ff8a855593417b12f2c3f4504fdc7a9c0dd4205/e2e/a.md.md", "", "", "a.md.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5a1d818eb87e27216135be2c7a065a94a404e74e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/6ff8a855593417b12f2c3f4504fdc7a9c0dd4205/e2e/a.md.md", "", "", "a.md.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5a1d818eb87e27216135be2c7a065a94a404e74e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf") | Out-Null

$de.Range("G2").Value = "2016-01-19 04:19:54"
$de.Range("G3").Value = "2016-01-19 04:19:54"

Write-Host "Handback report generated"
